# Update gh-pages output data for both "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# ---------- Sheet: 展览 ----------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 1226
$ws1.Range("F3").Value = 432

$ws1.Range("F5").Value = 12480
$ws1.Range("G5").Value = 60

$ws1.Range("F7").Value = 23
$ws1.Range("F8").Value = 25
$ws1.Range("F9").Value = 5

# Rows 10 and 11 swap content (with some updated values)
$ws1.Range("C10").Value = "苏州·I COME ACG动漫品牌博览会"
$ws1.Range("D10").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws1.Range("E10").Value = "2024.10.01 10:00-10.03 17:00"
$ws1.Range("F10").Value = 12357
$ws1.Range("G10").Value = 60
$ws1.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=87118"
$ws1.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg"

$ws1.Range("C11").Value = "苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场"
$ws1.Range("D11").Value = "苏州大道东688号 苏州国际博览中心"
$ws1.Range("E11").Value = "2024.10.01 09:00-10.01 17:00"
$ws1.Range("F11").Value = 231
$ws1.Range("G11").Value = 258
$ws1.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90770"
$ws1.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202408/reVAMzAd1723703136204.jpeg"

$ws1.Range("F12").Value = 4876
$ws1.Range("F13").Value = 4794
$ws1.Range("F14").Value = 149
$ws1.Range("F15").Value = 70
$ws1.Range("F17").Value = 102
$ws1.Range("F21").Value = 173
$ws1.Range("F22").Value = 77

# ---------- Sheet: 全部类型 ----------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 1226
$ws4.Range("F3").Value = 432

$ws4.Range("F7").Value = 12480
$ws4.Range("G7").Value = 60

$ws4.Range("F9").Value = 23
$ws4.Range("F10").Value = 25
$ws4.Range("F11").Value = 5

# Rows 12 and 13 swap content (with some updated values)
$ws4.Range("C12").Value = "苏州·I COME ACG动漫品牌博览会"
$ws4.Range("D12").Value = "金山南路288号木渎影视城F2 苏州广电国际会展中心"
$ws4.Range("E12").Value = "2024.10.01 10:00-10.03 17:00"
$ws4.Range("F12").Value = 12357
$ws4.Range("G12").Value = 60
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=87118"
$ws4.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg"

$ws4.Range("C13").Value = "苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场"
$ws4.Range("D13").Value = "苏州大道东688号 苏州国际博览中心"
$ws4.Range("E13").Value = "2024.10.01 09:00-10.01 17:00"
$ws4.Range("F13").Value = 231
$ws4.Range("G13").Value = 258
$ws4.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=90770"
$ws4.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202408/reVAMzAd1723703136204.jpeg"

$ws4.Range("F14").Value = 4876
$ws4.Range("F15").Value = 4794
$ws4.Range("F16").Value = 149
$ws4.Range("F17").Value = 70
$ws4.Range("F19").Value = 102
$ws4.Range("F23").Value = 173
$ws4.Range("F24").Value = 77
